# Commit: "A commit from my local computer"
#
# Adds a small yellow-highlighted reviewer note ("Make colors the same")
# as a text box on slide 1, placed above the existing pictures.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Create the text box (initial placement/size is irrelevant - refined below).
$tb = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)

# Exact geometry from the target OOXML, expressed in points (EMU / 12700):
#   a:off  x="3320143" y="381000"
#   a:ext cx="2235997" cy="369332"
$tb.Left   = 261.42858267716537
$tb.Top    = 30.0
$tb.Width  = 176.0627559055118
$tb.Height = 29.081259842519685

# Behaves like a freshly inserted PowerPoint text box: no word wrap, shape
# auto-sizes to fit its one line of text, and has no background fill.
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0

# Run text with a yellow (FFFF00) highlight.
$tr = $tb.TextFrame.TextRange
$tr.Text = "Make colors the same"
$tr.Font.Highlight.RGB = 65535
